$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Scattered single-cell corrections (outside the re-clustered block) ---
$ws.Range("I3").Value = 3
$ws.Range("E7").Value = 6
$ws.Range("J7").Value = 36
$ws.Range("J8").Value = 46
$ws.Range("E10").Value = 9
$ws.Range("K10").Value = 13
$ws.Range("K11").Value = 15
$ws.Range("I27").Value = 7

# --- Re-run clustering shifted rows 12-23: cluster 1 now ends at row 17
#     (Hurdman..Union) and cluster 2 begins at row 18 (Bayview..Wellesley),
#     with updated venue-count values from the new clustering pass. ---
# Row 12: Hurdman
$ws.Cells.Item(12, 1).Value = "Hurdman"
$ws.Cells.Item(12, 2).Value = "Ottawa"
$ws.Cells.Item(12, 3).Value = 45.412335300000002
$ws.Cells.Item(12, 4).Value = -75.664424862921507
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 11
$ws.Cells.Item(12, 7).Value = 15
$ws.Cells.Item(12, 8).Value = 3
$ws.Cells.Item(12, 9).Value = 14
$ws.Cells.Item(12, 10).Value = 16
$ws.Cells.Item(12, 11).Value = 15
$ws.Cells.Item(12, 12).Value = 25
$ws.Cells.Item(12, 13).Value = 1

# Row 13: Lees
$ws.Cells.Item(13, 1).Value = "Lees"
$ws.Cells.Item(13, 2).Value = "Ottawa"
$ws.Cells.Item(13, 3).Value = 45.41628455
$ws.Cells.Item(13, 4).Value = -75.6705328155996
$ws.Cells.Item(13, 5).Value = 4
$ws.Cells.Item(13, 6).Value = 15
$ws.Cells.Item(13, 7).Value = 18
$ws.Cells.Item(13, 8).Value = 4
$ws.Cells.Item(13, 9).Value = 12
$ws.Cells.Item(13, 10).Value = 17
$ws.Cells.Item(13, 11).Value = 12
$ws.Cells.Item(13, 12).Value = 19
$ws.Cells.Item(13, 13).Value = 1

# Row 14: Pimisi
$ws.Cells.Item(14, 1).Value = "Pimisi"
$ws.Cells.Item(14, 2).Value = "Ottawa"
$ws.Cells.Item(14, 3).Value = 45.413592350000002
$ws.Cells.Item(14, 4).Value = -75.713726352682102
$ws.Cells.Item(14, 5).Value = 15
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 22
$ws.Cells.Item(14, 8).Value = 1
$ws.Cells.Item(14, 9).Value = 13
$ws.Cells.Item(14, 10).Value = 19
$ws.Cells.Item(14, 11).Value = 18
$ws.Cells.Item(14, 12).Value = 10
$ws.Cells.Item(14, 13).Value = 1

# Row 15: Tremblay
$ws.Cells.Item(15, 1).Value = "Tremblay"
$ws.Cells.Item(15, 2).Value = "Ottawa"
$ws.Cells.Item(15, 3).Value = 45.416932799999998
$ws.Cells.Item(15, 4).Value = -75.653347895684107
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 20
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 5
$ws.Cells.Item(15, 10).Value = 16
$ws.Cells.Item(15, 11).Value = 25
$ws.Cells.Item(15, 12).Value = 30
$ws.Cells.Item(15, 13).Value = 1

# Row 16: Spadina
$ws.Cells.Item(16, 1).Value = "Spadina"
$ws.Cells.Item(16, 2).Value = "Toronto"
$ws.Cells.Item(16, 3).Value = 43.667234899999997
$ws.Cells.Item(16, 4).Value = -79.403686300000004
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 7
$ws.Cells.Item(16, 7).Value = 14
$ws.Cells.Item(16, 8).Value = 3
$ws.Cells.Item(16, 9).Value = 3
$ws.Cells.Item(16, 10).Value = 24
$ws.Cells.Item(16, 11).Value = 13
$ws.Cells.Item(16, 12).Value = 34
$ws.Cells.Item(16, 13).Value = 1

# Row 17: Union
$ws.Cells.Item(17, 1).Value = "Union"
$ws.Cells.Item(17, 2).Value = "Toronto"
$ws.Cells.Item(17, 3).Value = 43.644689999999997
$ws.Cells.Item(17, 4).Value = -79.379965688109493
$ws.Cells.Item(17, 5).Value = 4
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 26
$ws.Cells.Item(17, 8).Value = 4
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 19
$ws.Cells.Item(17, 11).Value = 13
$ws.Cells.Item(17, 12).Value = 33
$ws.Cells.Item(17, 13).Value = 1

# Row 18: Bayview
$ws.Cells.Item(18, 1).Value = "Bayview"
$ws.Cells.Item(18, 2).Value = "Ottawa"
$ws.Cells.Item(18, 3).Value = 45.409229850000003
$ws.Cells.Item(18, 4).Value = -75.722323334804599
$ws.Cells.Item(18, 5).Value = 6
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 16
$ws.Cells.Item(18, 8).Value = 7
$ws.Cells.Item(18, 9).Value = 6
$ws.Cells.Item(18, 10).Value = 26
$ws.Cells.Item(18, 11).Value = 32
$ws.Cells.Item(18, 12).Value = 7
$ws.Cells.Item(18, 13).Value = 2

# Row 19: Blair
$ws.Cells.Item(19, 1).Value = "Blair"
$ws.Cells.Item(19, 2).Value = "Ottawa"
$ws.Cells.Item(19, 3).Value = 45.431026250000002
$ws.Cells.Item(19, 4).Value = -75.608415591760505
$ws.Cells.Item(19, 5).Value = 1
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 27
$ws.Cells.Item(19, 8).Value = 1
$ws.Cells.Item(19, 9).Value = 4
$ws.Cells.Item(19, 10).Value = 20
$ws.Cells.Item(19, 11).Value = 35
$ws.Cells.Item(19, 12).Value = 10
$ws.Cells.Item(19, 13).Value = 2

# Row 20: Cyrville
$ws.Cells.Item(20, 1).Value = "Cyrville"
$ws.Cells.Item(20, 2).Value = "Ottawa"
$ws.Cells.Item(20, 3).Value = 45.422744950000002
$ws.Cells.Item(20, 4).Value = -75.626372478720796
$ws.Cells.Item(20, 5).Value = 2
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(20, 7).Value = 9
$ws.Cells.Item(20, 8).Value = 2
$ws.Cells.Item(20, 9).Value = 2
$ws.Cells.Item(20, 10).Value = 33
$ws.Cells.Item(20, 11).Value = 40
$ws.Cells.Item(20, 12).Value = 12
$ws.Cells.Item(20, 13).Value = 2

# Row 21: College
$ws.Cells.Item(21, 1).Value = "College"
$ws.Cells.Item(21, 2).Value = "Toronto"
$ws.Cells.Item(21, 3).Value = 43.660661699999999
$ws.Cells.Item(21, 4).Value = -79.382795200000004
$ws.Cells.Item(21, 5).Value = 4
$ws.Cells.Item(21, 6).Value = 2
$ws.Cells.Item(21, 7).Value = 24
$ws.Cells.Item(21, 8).Value = 3
$ws.Cells.Item(21, 9).Value = 4
$ws.Cells.Item(21, 10).Value = 28
$ws.Cells.Item(21, 11).Value = 28.999999999999901
$ws.Cells.Item(21, 12).Value = 7
$ws.Cells.Item(21, 13).Value = 2

# Row 22: Dundas
$ws.Cells.Item(22, 1).Value = "Dundas"
$ws.Cells.Item(22, 2).Value = "Toronto"
$ws.Cells.Item(22, 3).Value = 43.656536699999997
$ws.Cells.Item(22, 4).Value = -79.381022299999998
$ws.Cells.Item(22, 5).Value = 8
$ws.Cells.Item(22, 6).Value = 3
$ws.Cells.Item(22, 7).Value = 31
$ws.Cells.Item(22, 8).Value = 2
$ws.Cells.Item(22, 9).Value = 5
$ws.Cells.Item(22, 10).Value = 22
$ws.Cells.Item(22, 11).Value = 25
$ws.Cells.Item(22, 12).Value = 3
$ws.Cells.Item(22, 13).Value = 2

# Row 23: Wellesley
$ws.Cells.Item(23, 1).Value = "Wellesley"
$ws.Cells.Item(23, 2).Value = "Toronto"
$ws.Cells.Item(23, 3).Value = 43.665402999999998
$ws.Cells.Item(23, 4).Value = -79.383600099999995
$ws.Cells.Item(23, 5).Value = 4
$ws.Cells.Item(23, 6).Value = 1
$ws.Cells.Item(23, 7).Value = 25
$ws.Cells.Item(23, 8).Value = 6
$ws.Cells.Item(23, 9).Value = 3
$ws.Cells.Item(23, 10).Value = 16
$ws.Cells.Item(23, 11).Value = 38
$ws.Cells.Item(23, 12).Value = 7
$ws.Cells.Item(23, 13).Value = 2

# --- Restore view/selection state (scrolled back to top, new cursor cell) ---
$ws.Range("E21").Select()
